# close #125: Correct mandatory value columns assumptions in valores.xlsx
#
# The "1-2015" column (column C) was a mistaken/duplicate value column that
# should never have been included among the mandatory value columns — it
# simply repeated the "2-2015" figures. Remove it entirely so the sheet
# goes from 17 value columns (C:S) down to the correct 16 (C:R), with all
# subsequent columns shifting left to fill the gap.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("C").Delete()
